$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (A and B both to ~15.42578125; engine quantizes to 1/6 steps, closest achievable is 15.5)
$ws.Range("A1").EntireColumn.ColumnWidth = 14.666666666666668
$ws.Range("B1").EntireColumn.ColumnWidth = 14.666666666666668

# Update cell values
$ws.Range("A1").Value = -0.18419135961709543
$ws.Range("B1").Value = 0.1839390068102702
$ws.Range("A2").Value = -0.16183656714019357
$ws.Range("B2").Value = 0.16087237082842165
$ws.Range("A3").Value = -0.11116081581842785
$ws.Range("B3").Value = 0.11093457547303487
$ws.Range("A4").Value = -0.10293457557234476
$ws.Range("B4").Value = 0.1025452755606171
$ws.Range("A5").Value = -0.099545275618038609
$ws.Range("B5").Value = 0.09822653954633509
$ws.Range("A6").Value = -0.028297034817125422
$ws.Range("B6").Value = 0.028101069789721933
$ws.Range("A7").Value = -0.018101069934517433
$ws.Range("B7").Value = 0.018067409238734733
$ws.Range("A8").Value = -0.0080674093857413531
$ws.Range("B8").Value = 0.0080383166576827314
$ws.Range("A9").Value = -0.01239959987560546
$ws.Range("B9").Value = 0.012281093177163083
$ws.Range("A10").Value = 0.019380207292908125
$ws.Range("B10").Value = -0.019380124945303834
$ws.Range("A11").Value = 0.022380124862587003
$ws.Range("B11").Value = -0.02238596341740795
$ws.Range("A12").Value = -0.020864435766380307
$ws.Range("B12").Value = 0.020669837226262366
$ws.Range("A13").Value = -0.017169837318991021
$ws.Range("B13").Value = 0.017081571591431199
$ws.Range("A14").Value = -0.0090815717295296139
$ws.Range("B14").Value = 0.0090530500400811675
$ws.Range("A15").Value = -0.008053050111094251
$ws.Range("B15").Value = 0.008034653410452286
$ws.Range("A16").Value = -0.0060346534921489337
$ws.Range("B16").Value = 0.0060032556732165787
$ws.Range("A17").Value = -0.0040032557562756921
$ws.Range("B17").Value = 0.0039999998972435336
$ws.Range("A18").Value = -0.016102439730367735
$ws.Range("B18").Value = 0.016090988721103372
$ws.Range("A19").Value = -0.012090988762107013
$ws.Range("B19").Value = 0.012016121524238699
$ws.Range("A20").Value = -0.0080161215687404308
$ws.Range("B20").Value = 0.0080056186013290898
$ws.Range("A21").Value = -0.0040056186463299781
$ws.Range("B21").Value = 0.0039999999548285814
$ws.Range("A22").Value = -0.045711555064194442
$ws.Range("B22").Value = 0.045498751171056284
$ws.Range("A23").Value = -0.040498751238473574
$ws.Range("B23").Value = 0.040098819048469458
$ws.Range("A24").Value = -0.020098819268454804
$ws.Range("B24").Value = 0.019999999777065902
$ws.Range("A25").Value = -0.067929504795179341
$ws.Range("B25").Value = 0.067870650818727896
$ws.Range("A26").Value = -0.065370650891606985
$ws.Range("B26").Value = 0.065296537028867618
$ws.Range("A27").Value = -0.062796537104869543
$ws.Range("B27").Value = 0.062365733837574044
$ws.Range("A28").Value = -0.060365733922294496
$ws.Range("B28").Value = 0.060083601596630842
$ws.Range("A29").Value = -0.053083601739540853
$ws.Range("B29").Value = 0.053012235774057004
$ws.Range("A30").Value = 0.0069877635636994562
$ws.Range("B30").Value = -0.0070069132690759872
$ws.Range("A31").Value = -0.049025395156117924
$ws.Range("B31").Value = 0.048867500940964703
$ws.Range("A32").Value = -0.0040007950932796632
$ws.Range("B32").Value = 0.0039999998799356007
